$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each duty-roster name shifts down by one row (B2..B30 -> B3..B31),
# "高野怜央" is removed from the roster, and "小野文哉" is appended at the end.
$ws.Range("B2").Value = ""
$ws.Range("B3").Value = "志塚惇希"
$ws.Range("B6").Value = ""
$ws.Range("B7").Value = "白岩詩佑介"
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = "Nicholas Tristan Aryasatyo"
$ws.Range("B14").Value = ""
$ws.Range("B15").Value = "三神佳誠"
$ws.Range("B18").Value = ""
$ws.Range("B19").Value = "川田涼介"
$ws.Range("B22").Value = ""
$ws.Range("B23").Value = "兒島大志郎"
$ws.Range("B26").Value = ""
$ws.Range("B27").Value = "白岩詩佑介"
$ws.Range("B30").Value = ""
$ws.Range("B31").Value = "小野文哉"

# Update the selected cell to match the saved view state.
$ws.Range("D13").Select()
